# Update "想去人数" (F column) values on the 展览, 演出, and 全部类型 sheets
# to reflect newly scraped counts (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 194
$ws1.Range("F3").Value  = 5406
$ws1.Range("F7").Value  = 618
$ws1.Range("F8").Value  = 592
$ws1.Range("F9").Value  = 1058
$ws1.Range("F10").Value = 21
$ws1.Range("F11").Value = 1486
$ws1.Range("F12").Value = 4430
$ws1.Range("F17").Value = 3516
$ws1.Range("F18").Value = 178
$ws1.Range("F19").Value = 1110
$ws1.Range("F24").Value = 133
$ws1.Range("F26").Value = 144
$ws1.Range("F29").Value = 33
$ws1.Range("F30").Value = 58
$ws1.Range("F31").Value = 20
$ws1.Range("F33").Value = 32

# --- Sheet: 演出 ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 50

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 194
$ws4.Range("F3").Value  = 50
$ws4.Range("F4").Value  = 5406
$ws4.Range("F8").Value  = 618
$ws4.Range("F9").Value  = 592
$ws4.Range("F10").Value = 1058
$ws4.Range("F11").Value = 21
$ws4.Range("F12").Value = 1486
$ws4.Range("F13").Value = 4430
$ws4.Range("F18").Value = 3516
$ws4.Range("F19").Value = 178
$ws4.Range("F20").Value = 1110
$ws4.Range("F25").Value = 133
$ws4.Range("F27").Value = 144
$ws4.Range("F30").Value = 33
$ws4.Range("F31").Value = 58
$ws4.Range("F32").Value = 20
$ws4.Range("F34").Value = 32
